$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.213.00'
$ws.Range("E2").Value = '  +0.92%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.269.08'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.98'
$ws.Range("E5").Value = '  +0.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '96.91'
$ws.Range("E6").Value = '  +4.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.529'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  +1.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.64'
$ws.Range("E10").Value = '  +9.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0795'
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("E12").Value = '  -0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.64'
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.602.60'
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.42'
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.274.27'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.795'
$ws.Range("E17").Value = '  +2.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.135.60'
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.47'
$ws.Range("E19").Value = '  -1.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0909'
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.99'
$ws.Range("E21").Value = '  +1.18%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.87'
$ws.Range("E22").Value = '  +0.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.52'
$ws.Range("E23").Value = '  -2.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.58'
$ws.Range("E24").Value = '  +0.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.94'
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.75'
$ws.Range("E27").Value = '  -0.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.32'
$ws.Range("E28").Value = '  +6.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.53'
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.10'
$ws.Range("E30").Value = '  +1.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '159.57'
$ws.Range("E31").Value = '  -0.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.27'
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.18'
$ws.Range("E34").Value = '  +5.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0739'
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.31'
$ws.Range("E36").Value = '  +2.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.105'
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.37'
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("E39").Value = '  +2.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.114'
$ws.Range("E40").Value = '  -1.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.07'
$ws.Range("E41").Value = '  +4.00%  '
$ws.Range("E42").Value = '  +13.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.989.20'
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0285'
$ws.Range("E44").Value = '  +1.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.89'
$ws.Range("E45").Value = '  -4.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.93'
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.94'
$ws.Range("E47").Value = '  -3.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.15'
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.51'
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.11'
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '91.46'
$ws.Range("E51").Value = '  -0.03%  '
